$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.60186559259322
$ws.Range("C2").Value = 15.10079155969918
$ws.Range("D2").Value = 15.04923624578992
$ws.Range("E2").Value = 16.47344789375447
$ws.Range("G2").Value = 49.52032575154075
$ws.Range("H2").Value = 19.05053363868763
$ws.Range("I2").Value = 26.26998776359458
$ws.Range("J2").Value = 9.411846784306533
$ws.Range("B3").Value = 18.06459820391709
$ws.Range("C3").Value = 14.58921765728601
$ws.Range("D3").Value = 14.99556577656222
$ws.Range("E3").Value = 16.42078260341555
$ws.Range("G3").Value = 49.24180373085828
$ws.Range("H3").Value = 19.07581254945776
$ws.Range("I3").Value = 26.34656012163986
$ws.Range("J3").Value = 9.424415385093347
$ws.Range("B4").Value = 17.73078652705433
$ws.Range("C4").Value = 14.26981271441783
$ws.Range("D4").Value = 14.96626797246648
$ws.Range("E4").Value = 16.39247325198864
$ws.Range("G4").Value = 49.09010777741227
$ws.Range("H4").Value = 19.09648318718517
$ws.Range("I4").Value = 26.40141802202226
$ws.Range("J4").Value = 9.433707329963505
$ws.Range("B5").Value = 17.59398832533816
$ws.Range("C5").Value = 14.13853945566916
$ws.Range("D5").Value = 14.95525466792139
$ws.Range("E5").Value = 16.38195536308291
$ws.Range("G5").Value = 49.03317941890523
$ws.Range("H5").Value = 19.10619474344541
$ws.Range("I5").Value = 26.42573045349909
$ws.Range("J5").Value = 9.437889351207309
$ws.Range("B6").Value = 17.57123343897534
$ws.Range("C6").Value = 14.11668112873222
$ws.Range("D6").Value = 14.95348199818504
$ws.Range("E6").Value = 16.38027052768145
$ws.Range("G6").Value = 49.0240225108813
$ws.Range("H6").Value = 19.10788493116121
$ws.Range("I6").Value = 26.42988530058553
$ws.Range("J6").Value = 9.43860764233921
$ws.Range("B7").Value = 17.72894442466023
$ws.Range("C7").Value = 14.26804652271538
$ws.Range("D7").Value = 14.9661156870538
$ws.Range("E7").Value = 16.39232727395023
$ws.Range("G7").Value = 49.08932019269835
$ws.Range("H7").Value = 19.09660895459353
$ws.Range("I7").Value = 26.40173800239743
$ws.Range("J7").Value = 9.433762129622259
$ws.Range("B8").Value = 18.41755902037194
$ws.Range("C8").Value = 14.92562597230209
$ws.Range("D8").Value = 15.02997604783018
$ws.Range("E8").Value = 16.45445679007261
$ws.Range("G8").Value = 49.42030383091194
$ws.Range("H8").Value = 19.05817718460541
$ws.Range("I8").Value = 26.29475363681401
$ws.Range("J8").Value = 9.415853303505683
$ws.Range("B9").Value = 19.72777929127363
$ws.Range("C9").Value = 16.16440061315024
$ws.Range("D9").Value = 15.18388246715516
$ws.Range("E9").Value = 16.60794603795851
$ws.Range("G9").Value = 50.2207858006281
$ws.Range("H9").Value = 19.02394783947749
$ws.Range("I9").Value = 26.14781331876145
$ws.Range("J9").Value = 9.393250543403818
$ws.Range("B10").Value = 20.65491437900855
$ws.Range("C10").Value = 17.03319851627699
$ws.Range("D10").Value = 15.313923332214
$ws.Range("E10").Value = 16.73953426974677
$ws.Range("G10").Value = 50.89807844232674
$ws.Range("H10").Value = 19.02422263026427
$ws.Range("I10").Value = 26.0790070554962
$ws.Range("J10").Value = 9.384301363569172
$ws.Range("B11").Value = 21.06696265114226
$ws.Range("C11").Value = 17.41761756898315
$ws.Range("D11").Value = 15.37663452852016
$ws.Range("E11").Value = 16.80335895522751
$ws.Range("G11").Value = 51.22470844464529
$ws.Range("H11").Value = 19.02992949333699
$ws.Range("I11").Value = 26.05637303832021
$ws.Range("J11").Value = 9.381897166128802
$ws.Range("B12").Value = 21.22145003725028
$ws.Range("C12").Value = 17.56150252242429
$ws.Range("D12").Value = 15.40088013947292
$ws.Range("E12").Value = 16.82808556810081
$ws.Range("G12").Value = 51.35097521975063
$ws.Range("H12").Value = 19.03289711619911
$ws.Range("I12").Value = 26.04906096469591
$ws.Range("J12").Value = 9.381226646407818
$ws.Range("B13").Value = 21.18824942214702
$ws.Range("C13").Value = 17.53059126984154
$ws.Range("D13").Value = 15.39563646970608
$ws.Range("E13").Value = 16.82273566014089
$ws.Range("G13").Value = 51.32366813102319
$ws.Range("H13").Value = 19.0322220589129
$ws.Range("I13").Value = 26.05057956364599
$ws.Range("J13").Value = 9.381360382830312
$ws.Range("B14").Value = 21.07970421983601
$ws.Range("C14").Value = 17.4294895710313
$ws.Range("D14").Value = 15.37861932140882
$ws.Range("E14").Value = 16.80538212872059
$ws.Range("G14").Value = 51.23504531078065
$ws.Range("H14").Value = 19.03015745698455
$ws.Range("I14").Value = 26.05574616820229
$ws.Range("J14").Value = 9.381837193246199
$ws.Range("B15").Value = 21.01301156637203
$ws.Range("C15").Value = 17.36733858801727
$ws.Range("D15").Value = 15.36826031265942
$ws.Range("E15").Value = 16.79482481642839
$ws.Range("G15").Value = 51.18109450703624
$ws.Range("H15").Value = 19.02899797013985
$ws.Range("I15").Value = 26.05907517947335
$ws.Range("J15").Value = 9.382160500002298
$ws.Range("B16").Value = 20.62777749902213
$ws.Range("C16").Value = 17.0078467496946
$ws.Range("D16").Value = 15.309895456064
$ws.Range("E16").Value = 16.73544190660823
$ws.Range("G16").Value = 50.87709809616788
$ws.Range("H16").Value = 19.02396234298885
$ws.Range("I16").Value = 26.08066173420582
$ws.Range("J16").Value = 9.384492044144888
$ws.Range("B17").Value = 20.38885041373343
$ws.Range("C17").Value = 16.78444401560378
$ws.Range("D17").Value = 15.27499208918206
$ws.Range("E17").Value = 16.70001924641772
$ws.Range("G17").Value = 50.69529296573102
$ws.Range("H17").Value = 19.02230578279321
$ws.Range("I17").Value = 26.09613331202883
$ws.Range("J17").Value = 9.386349460480877
$ws.Range("B18").Value = 20.25051992471229
$ws.Range("C18").Value = 16.65493890035334
$ws.Range("D18").Value = 15.25525204050172
$ws.Range("E18").Value = 16.68001904390842
$ws.Range("G18").Value = 50.59247279840704
$ws.Range("H18").Value = 19.02187817437645
$ws.Range("I18").Value = 26.10584714860782
$ws.Range("J18").Value = 9.387574679658213
$ws.Range("B19").Value = 20.2035328148003
$ws.Range("C19").Value = 16.61092157376293
$ws.Range("D19").Value = 15.24862639193809
$ws.Range("E19").Value = 16.67331191428766
$ws.Range("G19").Value = 50.55796257842258
$ws.Range("H19").Value = 19.02182346743978
$ws.Range("I19").Value = 26.10927566333395
$ws.Range("J19").Value = 9.388016454812146
$ws.Range("B20").Value = 20.41437950587262
$ws.Range("C20").Value = 16.80833113835511
$ws.Range("D20").Value = 15.27867298588802
$ws.Range("E20").Value = 16.70375143157006
$ws.Range("G20").Value = 50.71446597585375
$ws.Range("H20").Value = 19.02242773977325
$ws.Range("I20").Value = 26.09440189292598
$ws.Range("J20").Value = 9.386135496836593
$ws.Range("B21").Value = 21.11162966747915
$ws.Range("C21").Value = 17.45923233633296
$ws.Range("D21").Value = 15.38360425591954
$ws.Range("E21").Value = 16.8104642532356
$ws.Range("G21").Value = 51.26100669710198
$ws.Range("H21").Value = 19.03074196323616
$ws.Range("I21").Value = 26.05419434597817
$ws.Range("J21").Value = 9.381690630442838
$ws.Range("B22").Value = 21.55825026582948
$ws.Range("C22").Value = 17.87475449537398
$ws.Range("D22").Value = 15.45508041556944
$ws.Range("E22").Value = 16.88345083583094
$ws.Range("G22").Value = 51.63319141288404
$ws.Range("H22").Value = 19.04087834354905
$ws.Range("I22").Value = 26.03525984544931
$ws.Range("J22").Value = 9.380184000807947
$ws.Range("B23").Value = 21.32075640357863
$ws.Range("C23").Value = 17.65392642841513
$ws.Range("D23").Value = 15.41667162989496
$ws.Range("E23").Value = 16.84420411503792
$ws.Range("G23").Value = 51.4332077647712
$ws.Range("H23").Value = 19.03503697392977
$ws.Range("I23").Value = 26.0446895935732
$ws.Range("J23").Value = 9.380860120042096
$ws.Range("B24").Value = 20.40284080726684
$ws.Range("C24").Value = 16.79753508619852
$ws.Range("D24").Value = 15.27700783497177
$ws.Range("E24").Value = 16.70206297359061
$ws.Range("G24").Value = 50.70579255019222
$ws.Range("H24").Value = 19.02237096878417
$ws.Range("I24").Value = 26.09518211741011
$ws.Range("J24").Value = 9.386231739657644
$ws.Range("B25").Value = 19.37882340619233
$ws.Range("C25").Value = 15.83585391665332
$ws.Range("D25").Value = 15.13922570968535
$ws.Range("E25").Value = 16.56307915718248
$ws.Range("G25").Value = 49.98831080081817
$ws.Range("H25").Value = 19.02876585531095
$ws.Range("I25").Value = 26.18074730538033
$ws.Range("J25").Value = 9.398022172077329
